$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Create the new sheet after sheet1 by copying (preserve formatting, cols, widths)
$ws1.Copy([Type]::Missing, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Edit_SYDEFAULT1"

# Remove the WO Picklist column (J) entirely - new sheet only has 9 columns (A:I)
$ws2.Columns.Item(10).Delete()

# Update row 2 values -> "No Override" for all data columns
$ws2.Range("A2:I2").Value = "No Override"

# Sheet1 view: no longer the tab-selected sheet; selection becomes A1:I2, scrolled to D1
$excel.ActiveWindow.ScrollColumn = 4
$ws1.Range("A1:I2").Select()

# Sheet2 view: becomes the active/tab-selected sheet; active cell H4
$ws2.Activate()
$ws2.Range("H4").Select()

Write-Host "done"
